$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 11:25"

# Filipinas (row 25) - refreshed totals
$ws.Range("B25").Value = 187249
$ws.Range("C25").Value = 4884
$ws.Range("D25").Value = 114921
$ws.Range("E25").Value = 69362
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 2966

# Indonesia (row 26) - refreshed totals
$ws.Range("B26").Value = 151498
$ws.Range("C26").Value = 2090
$ws.Range("D26").Value = 105198
$ws.Range("E26").Value = 39706
$ws.Range("G26").Value = 94
$ws.Range("H26").Value = 6594

# Ucrania overtakes Israel -> row 32 becomes Ucrania, row 33 becomes Israel
$ws.Range("A32").Value = "Ucrania"
$ws.Range("B32").Value = 102971
$ws.Range("C32").Value = 2328
$ws.Range("D32").Value = 51735
$ws.Range("E32").Value = 48992
$ws.Range("G32").Value = 37
$ws.Range("H32").Value = 2244

$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 100716
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 77785
$ws.Range("E33").Value = 22122
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 809

# Polonia (row 47) - refreshed totals
$ws.Range("B47").Value = 61181
$ws.Range("C47").Value = 900
$ws.Range("D47").Value = 41661
$ws.Range("E47").Value = 17569
$ws.Range("G47").Value = 13
$ws.Range("H47").Value = 1951

# Austria (row 71) - refreshed totals
$ws.Range("B71").Value = 25062
$ws.Range("C71").Value = 300
$ws.Range("D71").Value = 21406
$ws.Range("E71").Value = 2924
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 732

# Consejo Danes para los Refugiados (row 91) - refreshed totals
$ws.Range("B91").Value = 9811
$ws.Range("C91").Value = 9
$ws.Range("E91").Value = 640
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 251

# Hong Kong (row 112) - refreshed totals
$ws.Range("B112").Value = 4658
$ws.Range("C112").Value = 26
$ws.Range("D112").Value = 3974
$ws.Range("E112").Value = 608
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 76

# Sri Lanka (row 125) - refreshed totals
$ws.Range("D125").Value = 2798
$ws.Range("E125").Value = 132

# Lituania overtakes Eslovenia -> row 129 becomes Lituania, row 130 becomes Eslovenia
$ws.Range("A129").Value = "Lituania"
$ws.Range("B129").Value = 2594
$ws.Range("C129").Value = 30
$ws.Range("D129").Value = 1766
$ws.Range("E129").Value = 744
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 84

$ws.Range("A130").Value = "Eslovenia"
$ws.Range("B130").Value = 2574
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 2079
$ws.Range("E130").Value = 365
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 130
